# Worker List.xlsx update:
#  - Fill in newly-found scanner IDs for a few workers that didn't have one yet.
#  - Remove the "Changes" column (D), whose one-off annotations are no longer needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly discovered IDs for workers who previously had a blank ID cell.
# C5/C28 sit in "highlighted" rows, so line their ID cell up (left-aligned)
# with the rest of that row instead of leaving the centered blank style.
$ws.Range("C5").Value = "F76CB634"
$ws.Range("C5").HorizontalAlignment = -4131
$ws.Range("C6").Value = "4428B3A0"
$ws.Range("C28").Value = "B455EFD2"
$ws.Range("C28").HorizontalAlignment = -4131
$ws.Range("C32").Value = "0489D091"

# Drop the whole "Changes" column (D) -- header + notes.
$ws.Columns.Item(4).Delete()

# Refresh the remembered sort range now that the table is only A:C.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A1:A51"))
$ws.Sort.SetRange($ws.Range("A4:C51"))
$ws.Sort.Header = -4142
$ws.Sort.Apply()

# Match the final selection recorded in the saved file.
$ws.Range("C28").Select()
